$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 8 (2025) updates per diff:
# C8: 1171 -> 1178 (total_customers)
# D8: 193 (unchanged, returning_customers)
# E8: 978 -> 985 (new_customers)
# F8: unchanged (retention_rate)
# G8: 83.51836037574722 -> 83.61629881154499 (new_rate)
# H8: 16.48163962425278 -> 16.38370118845501 (returning_rate)

$ws.Range("C8").Value = 1178
$ws.Range("E8").Value = 985
$ws.Range("G8").Value = 83.61629881154499
$ws.Range("H8").Value = 16.38370118845501
